# edit.ps1
# Applies the "Horarios actualizados Linea 141 - 1275" update to
# horarios-141-2026-02-01.xlsx (3 sheets: LP1912, LP1912-215, 6203-6173).
#
# Summary of the update:
#  - "Ultima actualizacion" timestamp refreshed to 07:26:49 on all sheets.
#  - "Total filas" counters bumped on sheet1 (45->51) and sheet2 (9->10).
#  - Several existing rows refresh their scrape time (column A) and
#    recompute elapsed Minutos (column D) now that the scrape is newer;
#    a few rows also swap which "Linea" arrived first at a shared ETA.
#  - 6 new rows appended to sheet1 (51-56), 1 new row appended to sheet2 (15).

$wb = $excel.ActiveWorkbook

# ===== Sheet 1: LP1912 =====
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 07:26:49'
$ws1.Range("A3").Value = 'Total filas: 51'
$ws1.Range("A26").Value = '07:26:49'
$ws1.Range("D26").Value = 3
$ws1.Range("A28").Value = '07:26:49'
$ws1.Range("D28").Value = 9
$ws1.Range("A29").Value = '07:26:49'
$ws1.Range("D29").Value = 10
$ws1.Range("A31").Value = '07:26:49'
$ws1.Range("D31").Value = 17
$ws1.Range("A33").Value = '07:26:49'
$ws1.Range("D33").Value = 29
$ws1.Range("A35").Value = '07:26:49'
$ws1.Range("C35").Value = '16_SANTA ANA'
$ws1.Range("D35").Value = 34
$ws1.Range("A36").Value = '07:26:49'
$ws1.Range("C36").Value = '17_ROMERO'
$ws1.Range("D36").Value = 34
$ws1.Range("A38").Value = '07:26:49'
$ws1.Range("D38").Value = 40
$ws1.Range("A39").Value = '07:26:49'
$ws1.Range("D39").Value = 45
$ws1.Range("A41").Value = '07:26:49'
$ws1.Range("D41").Value = 47
$ws1.Range("C42").Value = '11_ETCHEVERRY'
$ws1.Range("C43").Value = '15_ABASTO'
$ws1.Range("A44").Value = '07:26:49'
$ws1.Range("C44").Value = '15_ABASTO'
$ws1.Range("D44").Value = 63
$ws1.Range("A45").Value = '07:26:49'
$ws1.Range("C45").Value = '11_ETCHEVERRY'
$ws1.Range("D45").Value = 63
$ws1.Range("A47").Value = '07:26:49'
$ws1.Range("D47").Value = 75
$ws1.Range("A48").Value = '07:26:49'
$ws1.Range("D48").Value = 77
$ws1.Range("A49").Value = '07:26:49'
$ws1.Range("B49").Value = '08:51'
$ws1.Range("D49").Value = 85
$ws1.Range("B50").Value = '08:52'
$ws1.Range("C50").Value = '23_HERNANDEZ'
$ws1.Range("D50").Value = 114
$ws1.Range("A51").Value = '07:26:49'
$ws1.Range("B51").Value = '08:53'
$ws1.Range("C51").Value = '215B_EL PATO'
$ws1.Range("D51").Value = 87
$ws1.Range("E51").Value = 'LP1912'
$ws1.Range("A52").Value = '07:26:49'
$ws1.Range("B52").Value = '08:57'
$ws1.Range("C52").Value = '215A_EL PATO'
$ws1.Range("D52").Value = 91
$ws1.Range("E52").Value = 'LP1912'
$ws1.Range("A53").Value = '07:26:49'
$ws1.Range("B53").Value = '09:06'
$ws1.Range("C53").Value = '16_SANTA ANA'
$ws1.Range("D53").Value = 100
$ws1.Range("E53").Value = 'LP1912'
$ws1.Range("A54").Value = '07:26:49'
$ws1.Range("B54").Value = '09:16'
$ws1.Range("C54").Value = '27_EL RETIRO'
$ws1.Range("D54").Value = 110
$ws1.Range("E54").Value = 'LP1912'
$ws1.Range("A55").Value = '07:26:49'
$ws1.Range("B55").Value = '09:17'
$ws1.Range("C55").Value = '14_ABASTO'
$ws1.Range("D55").Value = 111
$ws1.Range("E55").Value = 'LP1912'
$ws1.Range("A56").Value = '07:26:49'
$ws1.Range("B56").Value = '09:18'
$ws1.Range("C56").Value = '15X38_ABASTO'
$ws1.Range("D56").Value = 112
$ws1.Range("E56").Value = 'LP1912'

# ===== Sheet 2: LP1912-215 =====
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 07:26:49'
$ws2.Range("A3").Value = 'Total filas: 10'
$ws2.Range("A12").Value = '07:26:49'
$ws2.Range("D12").Value = 17
$ws2.Range("A13").Value = '07:26:49'
$ws2.Range("D13").Value = 77
$ws2.Range("A14").Value = '07:26:49'
$ws2.Range("D14").Value = 87
$ws2.Range("A15").Value = '07:26:49'
$ws2.Range("B15").Value = '08:57'
$ws2.Range("C15").Value = '215A_EL PATO'
$ws2.Range("D15").Value = 91
$ws2.Range("E15").Value = 'LP1912'

# ===== Sheet 3: 6203-6173 =====
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 07:26:49'
$ws3.Range("A6").Value = '07:26:49'
$ws3.Range("D6").Value = 16
$ws3.Range("A8").Value = '07:26:49'
$ws3.Range("D8").Value = 69
$ws3.Range("A9").Value = '07:26:49'
$ws3.Range("D9").Value = 84
